$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 45: C45 changes from text "en proceso" to numeric 1 formatted as percentage (task finished)
$ws.Range("C45").NumberFormat = "0%"
$ws.Range("C45").Value = 1

# Row 46: add responsible "Agustina" and mark as completed (100%)
$ws.Range("B46").Value = "Agustina"
$ws.Range("C46").NumberFormat = "0%"
$ws.Range("C46").Value = 1

# Row 47: add responsible "Agustina" and mark as completed (100%)
$ws.Range("B47").Value = "Agustina"
$ws.Range("C47").NumberFormat = "0%"
$ws.Range("C47").Value = 1

# Row 48: fix typo "rerportes" -> "reportes" in task description
$ws.Range("A48").Value = "reportes venta de productos (agrupar por producto/tipo, con codigo y combo/marca)"

# Update the sheet selection to match the new active cell/range
$ws.Range("B47:C47").Select() | Out-Null
